$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Park a copy of the original bordered date style (A2, style s=3) in a scratch cell
# so we can reapply it later after column formats are wiped.
$ws.Range("A2").Copy($ws.Range("Z1"))

# --- Remove the explicit column style (drops the <cols> block) and reset every
# cell in A:B back to the workbook default style (this also strips A1/B1 header styling).
$ws.Columns("A:B").ClearFormats()

# --- Build the "date, no border" style fresh on a scratch cell (now using the default font/fill).
$ws.Range("Z2").Value2 = 1
$ws.Range("Z2").NumberFormat = "mm-dd-yy"

$ws.Range("A2").Value2 = 40489
$ws.Range("B2").Value2 = 40492
$ws.Range("A3").Value2 = 40510
$ws.Range("B3").Value2 = 40522
$ws.Range("A4").Value2 = 40559
$ws.Range("B4").Value2 = 40562
$ws.Range("A5").Value2 = 40591
$ws.Range("B5").Value2 = 40595
$ws.Range("A6").Value2 = 40599
$ws.Range("B6").Value2 = 40623
$ws.Range("A7").Value2 = 40627
$ws.Range("B7").Value2 = 40638
$ws.Range("A8").Value2 = 40677
$ws.Range("B8").Value2 = 40684
$ws.Range("A9").Value2 = 40703
$ws.Range("B9").Value2 = 40706
$ws.Range("A10").Value2 = 40707
$ws.Range("B10").Value2 = 40762
$ws.Range("A11").Value2 = 40775
$ws.Range("B11").Value2 = 40865
$ws.Range("A12").Value2 = 40916
$ws.Range("B12").Value2 = 40957
$ws.Range("A13").Value2 = 40982
$ws.Range("B13").Value2 = 40993
$ws.Range("A14").Value2 = 41138
$ws.Range("B14").Value2 = 41140
$ws.Range("A15").Value2 = 41186
$ws.Range("B15").Value2 = 41208
$ws.Range("A16").Value2 = 41373
$ws.Range("B16").Value2 = 41380
$ws.Range("A17").Value2 = 41388
$ws.Range("B17").Value2 = 41397
$ws.Range("A18").Value2 = 41420
$ws.Range("B18").Value2 = 41461
$ws.Range("A19").Value2 = 41517
$ws.Range("B19").Value2 = 41549
$ws.Range("A20").Value2 = 41596
$ws.Range("B20").Value2 = 41597
$ws.Range("A21").Value2 = 41612
$ws.Range("B21").Value2 = 41615
$ws.Range("A22").Value2 = 41618
$ws.Range("B22").Value2 = 41626
$ws.Range("A23").Value2 = 41645
$ws.Range("B23").Value2 = 41694
$ws.Range("A24").Value2 = 41701
$ws.Range("B24").Value2 = 41739
$ws.Range("A25").Value2 = 41745
$ws.Range("B25").Value2 = 41778
$ws.Range("A26").Value2 = 41793
$ws.Range("B26").Value2 = 41815
$ws.Range("A27").Value2 = 41822
$ws.Range("B27").Value2 = 41917
$ws.Range("A28").Value2 = 41926
$ws.Range("B28").Value2 = 41945
$ws.Range("A29").Value2 = 41955
$ws.Range("B29").Value2 = 42018
$ws.Range("A30").Value2 = 42030
$ws.Range("B30").Value2 = 42035
$ws.Range("A31").Value2 = 42074
$ws.Range("B31").Value2 = 42156
$ws.Range("A32").Value2 = 42197
$ws.Range("B32").Value2 = 42234
$ws.Range("A33").Value2 = 42312
$ws.Range("B33").Value2 = 42332
$ws.Range("A34").Value2 = 42353
$ws.Range("B34").Value2 = 42364
$ws.Range("A35").Value2 = 42376
$ws.Range("B35").Value2 = 42384
$ws.Range("A36").Value2 = 42389
$ws.Range("B36").Value2 = 42400
$ws.Range("A37").Value2 = 42422
$ws.Range("B37").Value2 = 42434
$ws.Range("A38").Value2 = 42538
$ws.Range("B38").Value2 = 42543
$ws.Range("A39").Value2 = 42553
$ws.Range("B39").Value2 = 42614
$ws.Range("A40").Value2 = 42739
$ws.Range("B40").Value2 = 42746
$ws.Range("A41").Value2 = 42797
$ws.Range("B41").Value2 = 42818
$ws.Range("A42").Value2 = 42897
$ws.Range("B42").Value2 = 42932
$ws.Range("A43").Value2 = 42979
$ws.Range("B43").Value2 = 42992
$ws.Range("A44").Value2 = 43047
$ws.Range("B44").Value2 = 43051
$ws.Range("A45").Value2 = 43085
$ws.Range("B45").Value2 = 43099
$ws.Range("A46").Value2 = 43106
$ws.Range("B46").Value2 = 43136
$ws.Range("A47").Value2 = 43163
$ws.Range("B47").Value2 = 43196

# --- Paint the two style groups: rows with no border get the fresh Z2 style,
# rows 16-47 (old data + new trailing row) get the parked bordered style from Z1.
$ws.Range("Z2").Copy()
$ws.Range("A2:B15").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("A16:B47").PasteSpecial(-4122)

# --- Header row loses its bold/underline/fill/border styling -> default style (same trick: paste
# formats from an always-blank, always-default cell).
$ws.Range("Z3").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# --- Clean up the scratch column so it does not leak into the saved sheet / dimension.
$ws.Range("Z1:Z3").Clear()

# --- Window chrome: update the active selection + zoom to match the edited sheet state.
$ws.Range("B5").Select()
$win = $excel.ActiveWindow
$win.Zoom = 100

Write-Host "done"
